# Re-shuffle the "grid_cell" column (AG) on the "solar" sheet so that each
# elc_won/elc_spv distribution-process row points at a different CHE_n grid
# cell than before (same 25 labels, fully re-permuted, no row keeps its
# original value).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solar")

$newGridCell = @{
    4  = "CHE_24"
    5  = "CHE_8"
    6  = "CHE_5"
    7  = "CHE_3"
    8  = "CHE_10"
    9  = "CHE_22"
    10 = "CHE_12"
    11 = "CHE_20"
    12 = "CHE_1"
    13 = "CHE_6"
    14 = "CHE_0"
    15 = "CHE_13"
    16 = "CHE_9"
    17 = "CHE_21"
    18 = "CHE_4"
    19 = "CHE_2"
    20 = "CHE_14"
    21 = "CHE_18"
    22 = "CHE_17"
    23 = "CHE_19"
    24 = "CHE_23"
    25 = "CHE_11"
    26 = "CHE_15"
    27 = "CHE_25"
    28 = "CHE_7"
}

foreach ($row in $newGridCell.Keys) {
    $ws.Range("AG$row").Value = $newGridCell[$row]
}
